# Rename the existing (only) worksheet from "Sheet2" to "Sheet1"
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "Sheet1"

# Add a new worksheet right after Sheet1, named "Sheet2" (new sic_cat lookup table)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

$data = @(
    @("sic_cat", "cat"),
    @("A", "A"),
    @("B", "B-E"),
    @("C", "B-E"),
    @("D", "B-E"),
    @("E", "B-E"),
    @("F", "F"),
    @("G", "G-I"),
    @("H", "G-I"),
    @("I", "G-I"),
    @("J", "J"),
    @("K", "K"),
    @("L", "L"),
    @("M", "M-N"),
    @("N", "M-N"),
    @("O", "O-Q"),
    @("P", "O-Q"),
    @("Q", "O-Q"),
    @("R", "R-U"),
    @("S", "R-U"),
    @("T", "R-U"),
    @("U", "R-U")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws2.Cells.Item($row, 1).Value = $data[$i][0]
    $ws2.Cells.Item($row, 2).Value = $data[$i][1]
}

# Restore the selection on Sheet1 and set the new selection on Sheet2
$ws1.Range("C4").Select() | Out-Null
$ws2.Range("T13").Select() | Out-Null

# Sheet2 is the active tab
$ws2.Activate()
